# Daily attendance processing - 2026-01-07 10:38:21
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) "Recorded By" cells: flip the order "dnasr281@gmail.com, System"
#    -> "System, dnasr281@gmail.com" for the sessions that were
#    re-synced in this run.
# ------------------------------------------------------------------
$recordedByCells = @(
    "G2", "G3", "G4", "G5", "G6", "G7",
    "G16", "G17",
    "G22", "G23",
    "G37", "G38",
    "G43", "G44",
    "G58", "G59",
    "G64", "G65",
    "G79", "G80",
    "G85", "G86", "G87", "G88", "G89", "G90",
    "G99", "G100",
    "G105", "G106", "G107", "G108", "G109", "G110",
    "G119", "G120",
    "G125", "G126", "G127", "G128", "G129", "G130",
    "G139", "G140",
    "G145", "G146", "G147", "G148", "G149", "G150",
    "G159", "G160",
    "G165", "G166", "G167", "G168", "G169", "G170",
    "G179", "G180",
    "G185", "G186",
    "G200", "G201",
    "G206", "G207",
    "G221", "G222",
    "G227", "G228",
    "G242", "G243"
)

foreach ($addr in $recordedByCells) {
    $ws.Range($addr).Value = "System, dnasr281@gmail.com"
}

# ------------------------------------------------------------------
# 2) Six sessions (one per still-open group) crossed from "Pending"
#    into "Not Recorded" as their scheduled time passed. Re-stamp
#    their row formatting (green "Pending" -> pink "Not Recorded")
#    by pasting the format from an existing "Not Recorded" row, then
#    fix up the status label text.
# ------------------------------------------------------------------
$formatSource = $ws.Range("A18:I18")
$formatSource.Copy()

$notRecordedRows = @(31, 52, 73, 194, 215, 236)
foreach ($r in $notRecordedRows) {
    $target = $ws.Range("A" + $r + ":I" + $r)
    $target.PasteSpecial(-4122)
    $ws.Range("I" + $r).Value = "Not Recorded"
}

$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 3) Overall class statistics: six sessions moved from "Pending" to
#    "Missing", so totals shift by 6 in each direction.
# ------------------------------------------------------------------
$ws.Range("L7").Value = 54
$ws.Range("L8").Value = 72

# ------------------------------------------------------------------
# 4) Per-group statistics for the six affected groups: Missing +1,
#    Pending -1.
# ------------------------------------------------------------------
$groupStatRows = @(16, 17, 18, 24, 25, 26)
foreach ($r in $groupStatRows) {
    $ws.Range("P" + $r).Value = 4
    $ws.Range("Q" + $r).Value = 6
}
